$d = $word.ActiveDocument

# --- Edit 1: Education line -- replace "- Major (BS): ..." with "| BS in ..." ---
$oldEdu = " – Major (BS): Electrical & Computer Engineering | Minor: CS"
$newEdu = " | BS in Electrical & Computer Engineering | Minor: CS"
$found1 = $d.Content.Find.Execute($oldEdu, $true, $false, $false, $false, $false, $true, 1, $false, $newEdu, 2)
Write-Host "Education line replaced: $found1"

# --- Edit 2: GPA bullet -- "Current GPA of 3.59" -> "Current GPA: 3.62" ---
$oldGpa = "Current GPA of 3.59 – Dean’s List for Three Semesters"
$newGpa = "Current GPA: 3.62 – Dean’s List for Three Semesters"
$found2 = $d.Content.Find.Execute($oldGpa, $true, $false, $false, $false, $false, $true, 1, $false, $newGpa, 2)
Write-Host "GPA bullet replaced: $found2"

# --- Edit 3: Project bullet -- guitar hero MSP430 project -> pipe organ project ---
$oldProj = "Designed a guitar hero game on an MSP430 launchpad board. With this the user had to push the button they were instructed to, and if a certain amount of button presses were missed, the user would lose the game. The solo from Freebird played by Lynyrd Skynyrd was playing on a buzzer as well."
$newProj = "Developed a self-playing four-pipe organ on a team of four students. Designed and implemented the software to process MIDI data, control the stepper motors, and sync the solenoid for accurate timing."
$found3 = $d.Content.Find.Execute($oldProj, $true, $false, $false, $false, $false, $true, 1, $false, $newProj, 2)
Write-Host "Project bullet replaced: $found3"
